$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the target cells keep their original text ("inline string") type
# instead of being auto-converted to numbers by Excel when we assign
# numeric-looking strings.
$updates = @(
    @{Row=5; C="154"; D="408886.40"}
    @{Row=6; C="445"; D="1153510.82"}
    @{Row=7; C="199"; D="439041.00"}
    @{Row=8; C="811"; D="3083688.81"}
    @{Row=12; C="172"; D="533316.18"}
    @{Row=13; C="100"; D="242800.00"}
    @{Row=17; C="191"; D="431089.87"}
    @{Row=18; C="15"; D="56000.00"}
    @{Row=20; C="44"; D="121000.00"}
    @{Row=21; C="134"; D="376175.00"}
    @{Row=23; C="325"; D="1358916.35"}
    @{Row=24; C="6"; D="15000.00"}
    @{Row=25; C="9"; D="24024.69"}
    @{Row=26; C="39"; D="116572.00"}
    @{Row=27; C="14"; D="45500.00"}
    @{Row=28; C="16"; D="41450.00"}
    @{Row=31; C="56"; D="139310.00"}
    @{Row=33; C="31"; D="102000.00"}
    @{Row=35; C="169"; D="477408.00"}
    @{Row=37; C="380"; D="1517848.18"}
    @{Row=46; C="77"; D="176768.00"}
    @{Row=50; C="120"; D="301028.33"}
    @{Row=60; C="44"; D="194656.00"}
    @{Row=75; C="43"; D="128579.25"}
    @{Row=78; C="216"; D="601575.19"}
    @{Row=80; C="499"; D="2188776.03"}
    @{Row=85; C="38"; D="123669.00"}
    @{Row=88; C="74"; D="340136.08"}
    @{Row=120; C="34"; D="136158.69"}
    @{Row=122; C="254"; D="708508.00"}
    @{Row=123; C="127"; D="329012.45"}
    @{Row=124; C="512"; D="2325836.06"}
    @{Row=125; C="9"; D="37000.00"}
    @{Row=128; C="93"; D="280743.68"}
    @{Row=129; C="46"; D="178579.76"}
    @{Row=154; C="52"; D="177876.69"}
)

foreach ($u in $updates) {
    $cCell = $ws.Cells.Item($u.Row, 3)
    $dCell = $ws.Cells.Item($u.Row, 4)
    $cCell.NumberFormat = "@"
    $dCell.NumberFormat = "@"
    $cCell.Value = $u.C
    $dCell.Value = $u.D
}
